$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '96.620.93'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.664.49'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.82%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('E6').Value = '  +10.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '657.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.425'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.08'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.55%  '
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.661.69'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.32'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.49%  '
$ws.Range('E13').Value = '  +0.44%  '
$ws.Range('E14').Value = '  +4.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.346.96'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000272'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '96.307.94'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.675.25'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.18%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.84'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.39%  '
$ws.Range('E20').Value = '  +4.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.74'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('E22').Value = '  +1.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '530.65'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.12'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000205'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '102.21'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.06%  '
$ws.Range('E29').Value = '  +4.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.45'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.79%  '
$ws.Range('E31').Value = '  +0.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.93'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +17.39%  '
$ws.Range('E34').Value = '  +0.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '665.58'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.03%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '32.46'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.46%  '
$ws.Range('E38').Value = '  +4.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.85'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('E40').Value = '  +3.26%  '
$ws.Range('E41').Value = '  +0.71%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.54'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.74%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.957'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.81%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.71'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +17.08%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  +3.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.435'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +12.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.82'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.33'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.67'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.67'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.03%  '
